$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (logistic_embeddings)
$ws.Range("C5").Value = 0.5649999999999999
$ws.Range("D5").Value = 0.661
$ws.Range("E5").Value = 0.6889999999999999
$ws.Range("F5").Value = 0.722
$ws.Range("G5").Value = 0.621
$ws.Range("H5").Value = 0.622

# Row 7 (classical-best-embeddings -> classical-best-embed)
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.5649999999999999
$ws.Range("D7").Value = 0.661
$ws.Range("E7").Value = 0.6889999999999999
$ws.Range("F7").Value = 0.722
$ws.Range("H7").Value = 0.622

# Row 8 (BERT-base)
$ws.Range("C8").Value = 0.582
$ws.Range("D8").Value = 0.718
$ws.Range("E8").Value = 0.745
$ws.Range("F8").Value = 0.763
$ws.Range("G8").Value = 0.6909999999999999
$ws.Range("H8").Value = 0.699

# Row 9 (BERT-base-nli)
$ws.Range("B9").Value = 0.528
$ws.Range("C9").Value = 0.622
$ws.Range("D9").Value = 0.725
$ws.Range("E9").Value = 0.746
$ws.Range("F9").Value = 0.759
$ws.Range("G9").Value = 0.672
$ws.Range("H9").Value = 0.677
